# Commit: "pass targetTable in to CreateDestinationTable, load each
# worksheet if the xlsx file has more than one"
#
# This reflects a second worksheet ("primes") being added to the
# "alphabet" workbook (previously named "sheet.-name test"), so the
# bulk-insert tool has more than one sheet/table to exercise.

$wb = $excel.ActiveWorkbook

# --- Rename the original worksheet -----------------------------------
$alphabet = $wb.Worksheets.Item(1)
$alphabet.Name = "alphabet"

# Re-select the full table on the alphabet sheet (A1:B27) as it is in
# the saved workbook.
$alphabet.Range("A1:B27").Select() | Out-Null

# --- Add the new "primes" worksheet right after "alphabet" -----------
$primes = $wb.Worksheets.Add($null, $alphabet)
$primes.Name = "primes"

# Header row
$primes.Range("A1").Value = "Id"
$primes.Range("B1").Value = "PrimeNumber"

# Data rows: Id 1..11 next to the first 11 prime numbers
$primeNumbers = @(1, 2, 3, 5, 7, 11, 13, 17, 19, 23, 29)
for ($i = 0; $i -lt $primeNumbers.Length; $i++) {
    $row = $i + 2
    $primes.Cells.Item($row, 1).Value = $i + 1
    $primes.Cells.Item($row, 2).Value = $primeNumbers[$i]
}

# Widen column B ("PrimeNumber") to fit its header, matching the
# persisted column width of 14.7265625 as closely as this runtime's
# rounding allows.
$primes.Columns.Item(2).ColumnWidth = 13.893229166666666

# Selection on the new sheet, and make it the active tab.
$primes.Range("B4").Select() | Out-Null
$primes.Activate() | Out-Null
